# NCI Thesaurus and CIViC update
$wb = $excel.ActiveWorkbook

$compounds  = $wb.Worksheets.Item("compounds")
$biomarkers = $wb.Worksheets.Item("biomarkers")

# NCI Thesaurus source_version: 25.09e -> 25.10d
$compounds.Range("E3").Value = "25.10d"

# CIViC-related (Mitelman Database row) source_version stays v20250815
$biomarkers.Range("E3").Value = "v20250815"

# Update the selection/active state: "biomarkers" sheet's selection moves to A2
# and is no longer the active tab; "compounds" becomes the active tab while
# keeping its existing E3 selection.
$biomarkers.Range("A2").Select()
$compounds.Activate()
$compounds.Range("E3").Select()
